$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
Set-TextValue 2 4 "34.672.59"
Set-TextValue 2 5 "  -1.95%  "

# Row 3
Set-TextValue 3 4 "1.809.23"
Set-TextValue 3 5 "  -1.59%  "

# Row 4
Set-TextValue 4 5 "  +0.14%  "

# Row 5
Set-TextValue 5 4 "230.90"
Set-TextValue 5 5 "  +0.29%  "

# Row 6
Set-TextValue 6 4 "0.602"
Set-TextValue 6 5 "  -1.00%  "

# Row 7
Set-TextValue 7 5 "  +0.20%  "

# Row 8
Set-TextValue 8 4 "39.26"
Set-TextValue 8 5 "  -8.96%  "

# Row 9
Set-TextValue 9 4 "0.324"
Set-TextValue 9 5 "  +5.38%  "

# Row 10
Set-TextValue 10 5 "  -2.65%  "

# Row 11
Set-TextValue 11 4 "0.0991"
Set-TextValue 11 5 "  -1.84%  "

# Row 12
Set-TextValue 12 4 "2.071.44"
Set-TextValue 12 5 "  -1.65%  "

# Row 13
Set-TextValue 13 2 "WrappedEther"
Set-TextValue 13 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 13 4 "1.830.73"
Set-TextValue 13 5 "  -0.39%  "

# Row 14
Set-TextValue 14 2 "Chainlink"
Set-TextValue 14 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue 14 4 "11.22"
Set-TextValue 14 5 "  +0.28%  "

# Row 15
Set-TextValue 15 2 "Polygon"
Set-TextValue 15 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue 15 4 "0.667"
Set-TextValue 15 5 "  -0.16%  "

# Row 16
Set-TextValue 16 4 "4.59"
Set-TextValue 16 5 "  -1.52%  "

# Row 17
Set-TextValue 17 4 "34.682.89"
Set-TextValue 17 5 "  -1.86%  "

# Row 18
Set-TextValue 18 4 "69.41"
Set-TextValue 18 5 "  -0.69%  "

# Row 19
Set-TextValue 19 4 "0.0₃0783"
Set-TextValue 19 5 "  -1.76%  "

# Row 20
Set-TextValue 20 4 "240.58"
Set-TextValue 20 5 "  -1.51%  "

# Row 21
Set-TextValue 21 4 "11.85"
Set-TextValue 21 5 "  -1.32%  "

# Row 22
Set-TextValue 22 4 "4.68"
Set-TextValue 22 5 "  +0.01%  "

# Row 23
Set-TextValue 23 5 "  +0.24%  "

# Row 24
Set-TextValue 24 5 "  +2.22%  "

# Row 25
Set-TextValue 25 4 "171.11"
Set-TextValue 25 5 "  +1.23%  "

# Row 26
Set-TextValue 26 5 "  -2.02%  "

# Row 27
Set-TextValue 27 4 "17.17"
Set-TextValue 27 5 "  -2.96%  "

# Row 28
Set-TextValue 28 4 "0.121"
Set-TextValue 28 5 "  -0.87%  "

# Row 29
Set-TextValue 29 4 "1.52"
Set-TextValue 29 5 "  -0.80%  "

# Row 30
Set-TextValue 30 5 "  +0.16%  "

# Row 31
Set-TextValue 31 4 "4.04"
Set-TextValue 31 5 "  +3.07%  "

# Row 32
Set-TextValue 32 4 "0.0543"
Set-TextValue 32 5 "  -0.78%  "

# Row 33
Set-TextValue 33 4 "3.93"
Set-TextValue 33 5 "  -2.79%  "

# Row 34
Set-TextValue 34 4 "1.27"
Set-TextValue 34 5 "  +16.93%  "

# Row 35
Set-TextValue 35 5 "  -3.35%  "

# Row 36
Set-TextValue 36 4 "0.696"
Set-TextValue 36 5 "  +1.97%  "

# Row 37
Set-TextValue 37 4 "91.31"
Set-TextValue 37 5 "  -2.31%  "

# Row 38
Set-TextValue 38 5 "  +4.42%  "

# Row 39
Set-TextValue 39 4 "1.328.12"
Set-TextValue 39 5 "  -0.99%  "

# Row 40
Set-TextValue 40 5 "  -1.13%  "

# Row 41
Set-TextValue 41 5 "  +0.63%  "

# Row 42
Set-TextValue 42 4 "0.959"
Set-TextValue 42 5 "  -4.34%  "

# Row 43
Set-TextValue 43 4 "14.13"
Set-TextValue 43 5 "  -6.93%  "

# Row 44
Set-TextValue 44 5 "  -8.93%  "

# Row 45
Set-TextValue 45 5 "  -4.93%  "

# Row 46
Set-TextValue 46 4 "6.24"
Set-TextValue 46 5 "  +0.33%  "

# Row 47
Set-TextValue 47 5 "  -1.02%  "

# Row 48
Set-TextValue 48 4 "1.998.78"
Set-TextValue 48 5 "  -0.46%  "

# Row 49
Set-TextValue 49 5 "  +0.24%  "

# Row 50
Set-TextValue 50 5 "  +7.34%  "

# Row 51
Set-TextValue 51 4 "97.87"
Set-TextValue 51 5 "  -4.82%  "
